# "Our Clients" subsection - rename + restyle two client logos on slide 6.
# Shape.Left/Top/Width/Height are expressed in points (EMU / 12700); the
# literal values below are the closest points representations that
# round-trip to the target EMU offsets/extents from the canonical OOXML:
#   Rectangle 116: off x=13469538 (was 13498743), ext cx=736805 (was 678392)
#   Rectangle 117: off x=13304781 (was 13168524), ext cx=1066318 (was 1338829)
# (y/cy are untouched in both shapes.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Shape "Rectangle 116": Flickr -> Zillow
$flickr = $s.Shapes.Item(35)
$flickr.Left = 1060.593505859375
$flickr.Width = 58.016143798828125
$flickr.TextFrame.TextRange.Text = "Zillow"

# Shape "Rectangle 117": Belarusbank -> Sberbank
$sberbank = $s.Shapes.Item(36)
$sberbank.Left = 1047.62060546875
$sberbank.Width = 83.96205139160156
$sberbank.TextFrame.TextRange.Text = "Sberbank"
